# TC20_Canine_Filter_Breed-Dalmatian.xlsx
# "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The "startup" sheet's B2 cell held the Cypher query for the Cases tab.
# It previously returned a `Cohort` column; that column (and the MATCH/
# coalesce line that produced it) is removed here. B3/B4 (Samples/Files
# queries) are left untouched - their displayed text does not change,
# only their internal shared-string slot shifts as a natural side effect
# of de-duplicating strings after B2's text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Dalmatian']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesQuery

# The row shrank (one fewer wrapped line of text) - match the recalculated
# row heights for the three query rows.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 216
$ws.Rows.Item(4).RowHeight = 244.8

# Selection moved from D4 (with the view scrolled so row 4 was the top row)
# to B2, with the view back at the natural top-left.
[void]$ws.Activate()
[void]$ws.Range("B2").Select()
